{"js": "// Update the 25 \"three-digit \u00d7 one-digit\" equation prompts in the practice\n// table to a newly generated set of problems (output regenerated at c8c62b6).\n// Each old equation string is unique in the document, so a simple\n// search-and-replace per pair is unambiguous and safe.\nconst replacements = [\n  [\"444\u00d79=\", \"294\u00d75=\"],\n  [\"659\u00d78=\", \"473\u00d76=\"],\n  [\"916\u00d78=\", \"401\u00d73=\"],\n  [\"649\u00d73=\", \"344\u00d75=\"],\n  [\"223\u00d72=\", \"808\u00d77=\"],\n  [\"666\u00d72=\", \"572\u00d76=\"],\n  [\"729\u00d76=\", \"686\u00d76=\"],\n  [\"988\u00d72=\", \"336\u00d79=\"],\n  [\"578\u00d72=\", \"424\u00d76=\"],\n  [\"138\u00d78=\", \"342\u00d74=\"],\n  [\"293\u00d73=\", \"549\u00d72=\"],\n  [\"553\u00d73=\", \"748\u00d77=\"],\n  [\"853\u00d74=\", \"414\u00d74=\"],\n  [\"876\u00d74=\", \"158\u00d76=\"],\n  [\"596\u00d76=\", \"246\u00d74=\"],\n  [\"379\u00d75=\", \"815\u00d78=\"],\n  [\"310\u00d79=\", \"553\u00d74=\"],\n  [\"858\u00d77=\", \"800\u00d74=\"],\n  [\"620\u00d75=\", \"789\u00d75=\"],\n  [\"764\u00d77=\", \"532\u00d76=\"],\n  [\"750\u00d75=\", \"852\u00d76=\"],\n  [\"592\u00d73=\", \"949\u00d72=\"],\n  [\"950\u00d74=\", \"364\u00d78=\"],\n  [\"611\u00d76=\", \"408\u00d79=\"],\n  [\"562\u00d78=\", \"998\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"three-digit x one-digit\" equation prompts in the practice\n# table to a newly generated set of problems (output regenerated at c8c62b6).\n# Each old equation string is unique in the document, so Find/Replace per\n# pair is unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"444\u00d79=\", \"294\u00d75=\"),\n    @(\"659\u00d78=\", \"473\u00d76=\"),\n    @(\"916\u00d78=\", \"401\u00d73=\"),\n    @(\"649\u00d73=\", \"344\u00d75=\"),\n    @(\"223\u00d72=\", \"808\u00d77=\"),\n    @(\"666\u00d72=\", \"572\u00d76=\"),\n    @(\"729\u00d76=\", \"686\u00d76=\"),\n    @(\"988\u00d72=\", \"336\u00d79=\"),\n    @(\"578\u00d72=\", \"424\u00d76=\"),\n    @(\"138\u00d78=\", \"342\u00d74=\"),\n    @(\"293\u00d73=\", \"549\u00d72=\"),\n    @(\"553\u00d73=\", \"748\u00d77=\"),\n    @(\"853\u00d74=\", \"414\u00d74=\"),\n    @(\"876\u00d74=\", \"158\u00d76=\"),\n    @(\"596\u00d76=\", \"246\u00d74=\"),\n    @(\"379\u00d75=\", \"815\u00d78=\"),\n    @(\"310\u00d79=\", \"553\u00d74=\"),\n    @(\"858\u00d77=\", \"800\u00d74=\"),\n    @(\"620\u00d75=\", \"789\u00d75=\"),\n    @(\"764\u00d77=\", \"532\u00d76=\"),\n    @(\"750\u00d75=\", \"852\u00d76=\"),\n    @(\"592\u00d73=\", \"949\u00d72=\"),\n    @(\"950\u00d74=\", \"364\u00d78=\"),\n    @(\"611\u00d76=\", \"408\u00d79=\"),\n    @(\"562\u00d78=\", \"998\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
